$wb = $excel.ActiveWorkbook

# --- addScheduleForEmpTest sheet ---
$ws4 = $wb.Worksheets.Item("addScheduleForEmpTest")

# A2: "2022-09-11" -> "2022-09-09" (keep existing style/format; the leading
# apostrophe forces text entry so Excel doesn't reinterpret it as a real date)
$ws4.Range("A2").Value = "'2022-09-09"

# B2: "Manager Nine" -> "Hackshaw", and drop the one-off custom font formatting
# so the cell goes back to the default/general style.
$ws4.Range("B2").Value = "Hackshaw"
$ws4.Range("B2").ClearFormats()

# Selection on this (currently inactive) sheet moves from F20 to G13.
$ws4.Range("G13").Select()

# --- deleteScheduleForEmpTest sheet ---
$ws5 = $wb.Worksheets.Item("deleteScheduleForEmpTest")

# Content here is unchanged; select this sheet last so it stays the active tab,
# and widen the selection from A2 to the full row A2:XFD2.
$ws5.Range("A2:XFD2").Select()
